# Build site at 2022-09-26 16:07:08 UTC
#
# Rework of the "Docentes responsaveis" (responsible faculty) block on the
# LOM3086 sheet: the long free-text "Objetivos" / "Programa resumido" /
# "Programa" / "Bibliografia" paragraphs are dropped, the faculty roster
# (previously its own label-less rows) collapses onto the label rows that
# used to carry that removed text, and the remaining evaluation rows shift
# up to close the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the five rows that disappear from the middle of the block
#    (old rows 13-17), shrinking the sheet from 27 rows to 22.
$ws.Range("A13:A17").EntireRow.Delete()

# 2) Row 10 ("Objetivos:") now carries the first faculty entry instead of
#    the long objectives paragraph.
$ws.Range("B10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

# 3) Row 12 used to be the bare "Docentes responsáveis:" label (column A
#    only); it becomes "Programa resumido:" with the second faculty entry,
#    picking up the same 60pt custom row height used throughout this
#    block. B12/C12 are brand-new cells, so first clone the B/C number
#    format + alignment from a sibling row (row 14, still correctly
#    styled) before writing the value, matching the rest of the column.
$ws.Range("A12").Value = "Programa resumido:"
$ws.Rows.Item(12).RowHeight = 60
$ws.Range("B14").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("B12").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C12").Value = "5840897 - Clodoaldo Saron"

# 4) Row 14 ("Programa:") now carries the third faculty entry instead of
#    the long syllabus paragraph.
$ws.Range("B14").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C14").Value = "1033242 - Fábio Herbst Florenzano"

# 5) Row 17 ("Método:") now carries the fourth faculty entry instead of
#    the "Experimentos..." method text.
$ws.Range("B17").Value = "5840793 - Sérgio Schneider"
$ws.Range("C17").Value = "5840793 - Sérgio Schneider"

# 6) Everything below shifts up one slot: "Critério:" <- old "Método:"
#    text, "Norma de recuperação:" <- old "Critério:" text,
#    "Bibliografia:" <- old "Norma de recuperação:" text. The old
#    Bibliografia paragraph (the long reference list) is dropped entirely.
$ws.Range("B18").Value = "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento."
$ws.Range("C18").Value = "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento."

$ws.Range("B19").Value = "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."
$ws.Range("C19").Value = "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."

$ws.Range("B20").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
$ws.Range("C20").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
